$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose values are reset to 0 (forum view counters and derived totals/grades)
$cellsToZero = @(
    "E3", "G3", "I3", "J3", "G4", "I4", "J4", "H5", "I5", "J5",
    "F7", "I7", "J7", "B8", "C8", "F8", "I8", "J8", "C10", "D10",
    "H10", "I10", "J10", "D11", "F11", "I11", "J11", "C12", "D12", "E12",
    "F12", "I12", "J12", "C13", "I13", "J13", "B14", "D14", "E14", "I14",
    "J14", "C15", "E15", "F15", "I15", "J15", "D16", "I16", "J16", "F17",
    "I17", "J17", "D18", "F18", "I18", "J18", "B19", "C19", "F19", "G19",
    "I19", "J19", "D21", "I21", "J21", "D22", "E22", "F22", "I22", "J22",
    "B23", "C23", "D23", "E23", "F23", "G23", "I23", "J23", "C24", "D24",
    "I24", "J24", "B25", "C25", "I25", "J25", "C26", "D26", "E26", "G26",
    "I26", "J26", "F27", "I27", "J27", "C28", "D28", "E28", "F28", "I28",
    "J28", "D30", "F30", "G30", "I30", "J30", "C32", "D32", "E32", "H32",
    "I32", "J32", "B34", "C34", "I34", "J34", "B35", "C35", "G35", "I35",
    "J35", "F36", "I36", "J36", "D38", "E38", "G38", "I38", "J38", "B39",
    "E39", "G39", "H39", "I39", "J39", "D40", "E40", "F40", "G40", "I40",
    "J40", "H42", "I42", "J42", "G43", "I43", "J43", "B44", "D44", "E44",
    "F44", "G44", "I44", "J44", "B47", "C47", "F47", "I47", "J47", "H49",
    "I49", "J49", "C50", "G50", "I50", "J50"
)

foreach ($cellRef in $cellsToZero) {
    $ws.Range($cellRef).Value = 0
}

